$wb = $excel.ActiveWorkbook

# --- Content edits -------------------------------------------------------

# "Factory List": a new factory (blenderFactory) is being tracked, so add a
# row for it and make sure every factory row points at its own workbook
# file plus the standard chain-list / connections sheet names.
$ws1 = $wb.Worksheets.Item("Factory List")
$ws1.Activate()
$ws1.Range("A6").Value = "blenderFactory"

# The two cement factories shared a product name ("cement"); rename the
# duplicate's product on the 1990-Absolute scenario sheet so the industry
# balance can tell the two factories' outputs apart.
$ws2 = $wb.Worksheets.Item("1990-Absolute")
$ws2.Activate()
$ws2.Range("B4").Value = "same_cement"

# Back on "Factory List": fill in the Factory File / Chains / Connections
# columns for every factory (existing + new).
$ws1.Activate()
$ws1.Range("B2").Value = "excelData/cementFactory.xlsx"
$ws1.Range("B3").Value = "excelData/cementFactory.xlsx"
$ws1.Range("B4").Value = "excelData/otherCementFactory.xlsx"
$ws1.Range("B5").Value = "excelData/clinkerFactory.xlsx"
$ws1.Range("B6").Value = "excelData/blenderFactory.xlsx"

$ws1.Range("C2").Value = "Chain List"
$ws1.Range("C3").Value = "Chain List"
$ws1.Range("C4").Value = "Chain List"
$ws1.Range("C5").Value = "Chain List"
$ws1.Range("C6").Value = "Chain List"

$ws1.Range("D2").Value = "Connections"
$ws1.Range("D3").Value = "Connections"
$ws1.Range("D4").Value = "Connections"
$ws1.Range("D5").Value = "Connections"
$ws1.Range("D6").Value = "Connections"

# The "1990-2010" sheet actually reports the growth between the two years;
# rename it to make that explicit.
$ws6 = $wb.Worksheets.Item("1990-2010")
$ws6.Name = "1990-2010-Growth"

# --- Final selection / active-sheet state --------------------------------

$ws1.Range("C1").Select()
$ws2.Range("C33").Select()
$wb.Worksheets.Item("1990-Relative").Range("G38").Select()
$wb.Worksheets.Item("2010-Relative").Range("D7").Select()
$ws6.Range("J34").Select()

$wb.Worksheets.Item("2010-Relative").Activate()
